# "preparing runs for 2050" -- update Coupling Parameters sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# --- Year settings ---
$ws.Range("B2").Value = 2050   # Start Year
$ws.Range("B3").Value = 2060   # End Year
$ws.Range("B4").Value = 2050   # Power_plants_from_year

# --- Re-color B13:B16 to match the highlight used by B19:B29 (style index 7) ---
# Copy formatting from a cell that already carries that fill so the workbook
# reuses the existing style instead of creating a new one.
$ws.Range("B19").Copy()
$ws.Range("B13:B16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Investment / candidate capacity toggles ---
$ws.Range("B14").Value = $true   # realistic_candidate_capacities_tobe_installed
$ws.Range("B15").Value = $true   # realistic_candidate_capacities_for_future

# Row 17 ("testing_future_year") is shorter now that the note fits on fewer lines
$ws.Rows.Item(17).RowHeight = 16.5

# --- Fuel / price / demand settings ---
$ws.Range("B20").Value = 60       # start_tick_fuel_trends
$ws.Range("B21").Value = $false   # yearly_CO2_prices
$ws.Range("B22").Value = $true    # fix_fuel_prices_to_year
$ws.Range("B23").Value = 2050     # fix_price_year
$ws.Range("B24").Value = $true    # fix_demand_to_initial_year
$ws.Range("B25").Value = $true    # fix_profiles_to_initial_year

# --- Decommissioning / investment toggles ---
$ws.Range("B26").Value = $false   # decommission_from_input
$ws.Range("B27").Value = $false   # targetinvestment_per_year

# --- Selection left where the user was last working ---
$ws.Activate()
$ws.Range("B15").Select()
